$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 0.74517412566490038
$ws.Range("Z1").Value = 0.79160062094723305
$ws.Range("BP1").Value = 0.97221262956888943
$ws.Range("D2").Value = 0.74716357614557327
$ws.Range("B3").Value = 0.72870219841411332
$ws.Range("D3").Value = 0.92425409223929522
$ws.Range("J3").Value = 0.94799492440595556
$ws.Range("AP4").Value = 0.57255281313751039
$ws.Range("BH4").Value = 0.93514329700205234
$ws.Range("C5").Value = 0.63084609122173374
$ws.Range("G5").Value = 0.79921329165262267
$ws.Range("E6").Value = 0.53538647999659017
$ws.Range("H6").Value = 0.91252567426253206
$ws.Range("F7").Value = 0.93282693817278539
$ws.Range("H7").Value = 0.73217341061473529
$ws.Range("I7").Value = 0.87094597167919496
$ws.Range("I8").Value = 0.84936904706899274
$ws.Range("J9").Value = 0.76025383630080479
$ws.Range("H10").Value = 0.97487634826343084
$ws.Range("K10").Value = 0.65670297376128672
$ws.Range("I11").Value = 0.98188132259466476
$ws.Range("P11").Value = 0.71898207233683165
$ws.Range("K12").Value = 0.94566785698247202
$ws.Range("M12").Value = 0.58732532014825378
$ws.Range("K13").Value = 0.71891953368444139
$ws.Range("N13").Value = 0.82692955199694174
$ws.Range("L14").Value = 0.8688105738038423
$ws.Range("O14").Value = 0.93687401715473295
$ws.Range("M15").Value = 0.92807088170427365
$ws.Range("Q15").Value = 0.69270593488628207
$ws.Range("AO15").Value = 0.80930887408196006
$ws.Range("R16").Value = 0.95150161638224418
$ws.Range("R17").Value = 0.99384327106431369
$ws.Range("S17").Value = 0.71270844549094248
$ws.Range("AY17").Value = 0.82335502307592068
$ws.Range("S18").Value = 0.8595457389399801
$ws.Range("T18").Value = 0.88347395092729264
$ws.Range("T19").Value = 0.95539301603857463
$ws.Range("BM19").Value = 0.78250470388061499
$ws.Range("V20").Value = 0.65988182486516056
$ws.Range("S21").Value = 0.69002542564971014
$ws.Range("T21").Value = 0.69369858368416459
$ws.Range("V21").Value = 0.52343940509220843
$ws.Range("W22").Value = 0.89374318835386224
$ws.Range("X22").Value = 0.88627742863540782
$ws.Range("U23").Value = 0.98550151236661765
$ws.Range("W24").Value = 0.68192754763818386
$ws.Range("Y24").Value = 0.97730900351194838
$ws.Range("AZ25").Value = 0.95801600518019825
$ws.Range("X26").Value = 0.91053410419586778
$ws.Range("AA26").Value = 0.89971256572979241
$ws.Range("T27").Value = 0.96912307430330713
$ws.Range("Y27").Value = 0.86714673150906763
$ws.Range("AC27").Value = 0.96443877478645534
$ws.Range("A28").Value = 0.97358053060800998
$ws.Range("Z28").Value = 0.90180320552274007
$ws.Range("BE28").Value = 0.9031587534093809
$ws.Range("W29").Value = 0.77829050708493719
$ws.Range("AE29").Value = 0.94219154141480765
$ws.Range("AC30").Value = 0.81255792454177067
$ws.Range("AD31").Value = 0.79137349895058395
$ws.Range("AD32").Value = 0.68692741503324439
$ws.Range("AE32").Value = 0.84149216230192347
$ws.Range("AH32").Value = 0.92874719741339984
$ws.Range("AE33").Value = 0.99745486318955168
$ws.Range("AF33").Value = 0.65820840508176581
$ws.Range("AT33").Value = 0.85176019577501239
$ws.Range("AG34").Value = 0.93750851071182795
$ws.Range("AI34").Value = 0.77194640960916061
$ws.Range("AG35").Value = 0.63677442954276287
$ws.Range("BC35").Value = 0.86603053967031574
$ws.Range("AB36").Value = 0.83640828589163496
$ws.Range("AH36").Value = 0.78396321785564327
$ws.Range("AI36").Value = 0.99334746584120248
$ws.Range("AI37").Value = 0.75396857715357368
$ws.Range("AJ37").Value = 0.70287738644150877
$ws.Range("AM37").Value = 0.93908838142989293
$ws.Range("AM38").Value = 0.63449561341506222
$ws.Range("AN39").Value = 0.85795240272433204
$ws.Range("AX39").Value = 0.80610484700144158
$ws.Range("BD39").Value = 0.74721651277647294
$ws.Range("X40").Value = 0.89837551467614318
$ws.Range("AL40").Value = 0.87730412158299187
$ws.Range("BP40").Value = 0.99322242435766928
$ws.Range("AP41").Value = 0.84459447355042105
$ws.Range("AQ41").Value = 0.87136196409275457
$ws.Range("BI41").Value = 0.82652363276552221
$ws.Range("AN42").Value = 0.82423825245610427
$ws.Range("AP43").Value = 0.88791917267465081
$ws.Range("AS43").Value = 0.94201335176855261
$ws.Range("N44").Value = 0.8264292644891269
$ws.Range("AP44").Value = 0.9532599638059096
$ws.Range("AQ44").Value = 0.88233850034164862
$ws.Range("AR45").Value = 0.88610129112793912
$ws.Range("BI45").Value = 0.71007278660716922
$ws.Range("AR46").Value = 0.92228649470742063
$ws.Range("AO47").Value = 0.92035654689419188
$ws.Range("AS47").Value = 0.85777643651276403
$ws.Range("AH48").Value = 0.61785391721022664
$ws.Range("AT48").Value = 0.99332466812749454
$ws.Range("AU48").Value = 0.77567595148481172
$ws.Range("AN50").Value = 0.9203811415235037
$ws.Range("AW50").Value = 0.87340439630317346
$ws.Range("AW51").Value = 0.91142260697652144
$ws.Range("AZ51").Value = 0.69084962362201208
$ws.Range("BA52").Value = 0.63724802155530647
$ws.Range("AK53").Value = 0.93533114568726261
$ws.Range("BB53").Value = 0.71042236050172014
$ws.Range("BC53").Value = 0.73628807520570128
$ws.Range("AZ54").Value = 0.9122543244043736
$ws.Range("BD54").Value = 0.73157400780060522
$ws.Range("BB55").Value = 0.91477844054020774
$ws.Range("AI56").Value = 0.79490331411098225
$ws.Range("BC56").Value = 0.89979532315643196
$ws.Range("AK57").Value = 0.99660686793270004
$ws.Range("BF57").Value = 0.68524592568726361
$ws.Range("BD58").Value = 0.78408790701891884
$ws.Range("BG58").Value = 0.7965663843501718
$ws.Range("AH59").Value = 0.85578402072154147
$ws.Range("BI59").Value = 0.93322542721061241
$ws.Range("BF60").Value = 0.93981412743043113
$ws.Range("BJ60").Value = 0.74545247893996902
$ws.Range("BH61").Value = 0.69011446084763073
$ws.Range("BJ61").Value = 0.79328688385711033
$ws.Range("BL62").Value = 0.69249601109449221
$ws.Range("BJ63").Value = 0.87936274469525877
$ws.Range("BM63").Value = 0.8869152518407869
$ws.Range("BI64").Value = 0.83797133269064172
$ws.Range("BK64").Value = 0.62122487049200492
$ws.Range("BM64").Value = 0.80115490405244572
$ws.Range("BN65").Value = 0.87156483942753504
$ws.Range("BO66").Value = 0.82776898766884321
$ws.Range("BP66").Value = 0.79041348128841138
$ws.Range("A67").Value = 0.98469305646169958
$ws.Range("N67").Value = 0.92603853951418924
$ws.Range("AV67").Value = 0.86755982046729796
